$wb = $excel.ActiveWorkbook

# Update the "Advisor" description cells on the education sheet to include
# markdown links to the advisors' personal/lab websites.
$eduSheet = $wb.Worksheets.Item("education")
$eduSheet.Range("G2").Value = "Advisor: [Neal Kingston, Ph.D.](https://nealkingston.ku.edu/)"
$eduSheet.Range("G3").Value = "Advisor: [Evangelia G. Chrysikou, Ph.D.](https://www.chrysikoulab.com/)"

# Make "education" the active sheet/tab with G4 selected, matching the
# author's final view state captured in the commit.
$eduSheet.Activate()
$eduSheet.Range("G4").Select()
